# Generate Report for Handoff
# - Flip the "In Translation" status to "Ready for handoff" everywhere it
#   appears (Overview + per-locale sheets), and refresh the associated
#   handoff timestamps.
# - Widen the Status column(s) to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-locale status + last handoff-generation time ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 06:39:54"

# --- zh-cn sheet: status + latest handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 06:39:49"

# --- de-de sheet: status + latest handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 06:39:54"

# --- Widen the Status columns so the longer "Ready for handoff" text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.3   # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3       # column C (Status)
